$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 240, pushing existing data down by 2 rows
$ws.Rows("240:241").Insert()

# Fill in new row 240
$ws.Cells.Item(240, 1).Value = 8
$ws.Cells.Item(240, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(240, 3).Value = "Coquimbo"
$ws.Cells.Item(240, 4).Value = "2022-06-14"
$ws.Cells.Item(240, 5).Value = 4
$ws.Cells.Item(240, 6).Value = 100112032
$ws.Cells.Item(240, 7).Value = "Zapallo italiano"
$ws.Cells.Item(240, 8).Value = "Bola 8"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 600
$ws.Cells.Item(240, 11).Value = 10500
$ws.Cells.Item(240, 12).Value = 11000
$ws.Cells.Item(240, 13).Value = 10750
$ws.Cells.Item(240, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(240, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(240, 16).Value = 215
$ws.Cells.Item(240, 17).Value = 50
$ws.Cells.Item(240, 18).Value = "Hortaliza"

# Fill in new row 241
$ws.Cells.Item(241, 1).Value = 8
$ws.Cells.Item(241, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(241, 3).Value = "Coquimbo"
$ws.Cells.Item(241, 4).Value = "2022-06-14"
$ws.Cells.Item(241, 5).Value = 4
$ws.Cells.Item(241, 6).Value = 100112032
$ws.Cells.Item(241, 7).Value = "Zapallo italiano"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 500
$ws.Cells.Item(241, 11).Value = 11000
$ws.Cells.Item(241, 12).Value = 12000
$ws.Cells.Item(241, 13).Value = 11500
$ws.Cells.Item(241, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(241, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(241, 16).Value = 230
$ws.Cells.Item(241, 17).Value = 50
$ws.Cells.Item(241, 18).Value = "Hortaliza"
